$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append new row 9 to the Logs sheet
$ws.Range("A9").Value = "Retour status"
$ws.Range("B9").Value = "mailmind.test@zohomail.eu"
$ws.Range("D9").Value = "Retour / Terugbetaling"
$ws.Range("F9").Value = "2025-08-26 21:08:11"
$ws.Range("G9").Value = "Nee"
$ws.Range("H9").Value = "Ja"
$ws.Range("I9").Value = "Nee"
$ws.Range("J9").Value = "Nee"

# Extend the conditional formatting ranges to include the new row
$ws.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D9"))
$ws.Range("G2:G8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G9"))
$ws.Range("H2:H8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H9"))
$ws.Range("I2:I8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I9"))
$ws.Range("J2:J8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J9"))

# Update the Dashboard summary count for "Retour / Terugbetaling"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B3").Value = 4
